# Generate Report for Handoff
# Update the Priority and Latest Handoff Datetime columns for the rows
# that are in "Ready for handoff" status on both the zh-cn and de-de
# localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-17 04:28:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-17 04:28:31"
